$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update bound values (new values for upper and lower bound)
$ws.Range("C2").Value = 0.1

$ws.Range("B5").Value = 0.25
$ws.Range("C5").Value = 0.05

$ws.Range("B7").Value = 0.3
$ws.Range("C7").Value = 0.1

$ws.Range("B8").Value = 0.2

$ws.Range("B12").Value = 0.25
$ws.Range("C12").Value = 0.05

$ws.Range("B13").Value = 0.25
$ws.Range("C13").Value = 0.05

$ws.Range("B14").Value = 0.25

$ws.Range("B16").Value = 0.25

$ws.Range("B17").Value = 0.25

$ws.Range("B18").Value = 0.25
$ws.Range("C18").Value = 0.05

$ws.Range("B19").Value = 0.25
$ws.Range("C19").Value = 0.05

# Update selection to B10
$ws.Range("B10").Select()
